# Update timestamps in the handback status report, as if a new report
# generation run had completed slightly later than the previous one.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# Overview sheet: "Latest HO Xliff Generate Date" for first row (G2)
$overview.Range("G2").Value = "2016-08-26 01:04:10"

# zh-cn sheet: "Correspond Handoff Datetime" (H2) and
# "Correspond Handback DateTime" (K2) for first row
$zhcn.Range("H2").Value = "2016-08-26 01:04:01"
$zhcn.Range("K2").Value = "2016-08-26 01:04:39"

# de-de sheet: "Correspond Handoff Datetime" (H2) shares the same value as
# the Overview sheet's "Latest HO Xliff Generate Date", and
# "Correspond Handback DateTime" (K2) for first row
$dede.Range("H2").Value = "2016-08-26 01:04:10"
$dede.Range("K2").Value = "2016-08-26 01:04:45"
